$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the four new "Line" connect/disconnect rows to the Snippets table.
$lo = $ws.ListObjects.Item("Snippets")
$lo.Resize($ws.Range("A1:D173"))

$rows = @(
    @("Line", "connectBeginShape",    "excel-shape-lines", "connectStraightLine"),
    @("Line", "connectEndShape",      "excel-shape-lines", "connectStraightLine"),
    @("Line", "disconnectBeginShape", "excel-shape-lines", "disconnectStraightLine"),
    @("Line", "disconnectEndShape",   "excel-shape-lines", "disconnectStraightLine")
)

# Populate column D first, then A, then B, then C -- matching the shared-string
# insertion order recorded by the original authoring session (new unique
# strings are appended to sharedStrings.xml in first-use order).
$r = 170
foreach ($row in $rows) {
    $ws.Range("D$r").Value = $row[3]
    $r = $r + 1
}

$r = 170
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row[0]
    $r = $r + 1
}

$r = 170
foreach ($row in $rows) {
    $ws.Range("B$r").Value = $row[1]
    $r = $r + 1
}

$r = 170
foreach ($row in $rows) {
    $ws.Range("C$r").Value = $row[2]
    $r = $r + 1
}

# Columns A, B and D on the new rows carry the same (banded/general) cell
# style used elsewhere in the sheet; re-applying "General" stamps that xf.
$r = 170
foreach ($row in $rows) {
    $ws.Range("A$r").NumberFormat = "General"
    $ws.Range("B$r").NumberFormat = "General"
    $ws.Range("D$r").NumberFormat = "General"
    $r = $r + 1
}

$ws.Range("B173").Select()
